$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Core content change (commit: "add language id for property name")
# Row 9 held the Chinese-language descriptions for each property column.
# Replace them with language-id keys ("LPID_" + the machine name already
# present in row 1) so the sheet can be localized.  Overwriting these cells
# drops the old Chinese shared strings (they become unreferenced) and
# appends the new LPID_* strings to the shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "LPID_SUCKBLOOD"
$ws.Range("C9").Value = "LPID_REFLECTDAMAGE"
$ws.Range("D9").Value = "LPID_CRITICAL"
$ws.Range("E9").Value = "LPID_MAXHP"
$ws.Range("F9").Value = "LPID_MAXMP"
$ws.Range("G9").Value = "LPID_MAXSP"
$ws.Range("H9").Value = "LPID_HPREGEN"
$ws.Range("I9").Value = "LPID_SPREGEN"
$ws.Range("J9").Value = "LPID_MPREGEN"
$ws.Range("K9").Value = "LPID_ATK_VALUE"
$ws.Range("L9").Value = "LPID_DEF_VALUE"
$ws.Range("M9").Value = "LPID_MOVE_SPEED"
$ws.Range("N9").Value = "LPID_ATK_SPEED"
$ws.Range("O9").Value = "LPID_ATK_FIRE"
$ws.Range("P9").Value = "LPID_ATK_LIGHT"
$ws.Range("Q9").Value = "LPID_ATK_WIND"
$ws.Range("R9").Value = "LPID_ATK_ICE"
$ws.Range("S9").Value = "LPID_ATK_POISON"
$ws.Range("T9").Value = "LPID_DEF_FIRE"
$ws.Range("U9").Value = "LPID_DEF_LIGHT"
$ws.Range("V9").Value = "LPID_DEF_WIND"
$ws.Range("W9").Value = "LPID_DEF_ICE"
$ws.Range("X9").Value = "LPID_DEF_POISON"
$ws.Range("Y9").Value = "LPID_DIZZY_GATE"
$ws.Range("Z9").Value = "LPID_MOVE_GATE"
$ws.Range("AA9").Value = "LPID_SKILL_GATE"
$ws.Range("AB9").Value = "LPID_PHYSICAL_GATE"
$ws.Range("AC9").Value = "LPID_MAGIC_GATE"
$ws.Range("AD9").Value = "LPID_BUFF_GATE"

# ---------------------------------------------------------------------------
# Column width tweaks that came along with the edit (column B got its own
# width once it stopped sharing the merged Chinese-description formatting,
# and the newly meaningful AA:AD columns got explicit widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 15.571428571428571
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 8.571428571428571
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 11
$ws.Columns.Item(10).ColumnWidth = 14.571428571428571
$ws.Range("L1:Z1").EntireColumn.ColumnWidth = 15.857142857142858
$ws.Columns.Item(27).ColumnWidth = 9.571428571428571
$ws.Columns.Item(28).ColumnWidth = 12.285714285714286
$ws.Range("AC1:AD1").EntireColumn.ColumnWidth = 11.428571428571429

# ---------------------------------------------------------------------------
# View-state: selection moved to the last header cell (AD9) after the edit.
# ---------------------------------------------------------------------------
$ws.Range("AD9").Select()
